# Applies the "Made small changes to final section of poster" edit:
#  - TextBox 20 (shape 4): split the closing paragraph in two, inserting a
#    new blank paragraph and a new paragraph about accessibility features,
#    and drop the old "fully secure back-end" / chatbot sentence from the
#    first paragraph (it now opens the new paragraph instead).
#  - TextBox 23 (shape 5): merge the three runs of the last paragraph into
#    one and fix the "guidlines" typo to "guidelines".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# Shape 4 ("TextBox 20") - closing paragraph of the "Conclusion" text
# ---------------------------------------------------------------
$shape20 = $s.Shapes.Item("TextBox 20")
$tr20 = $shape20.TextFrame.TextRange

$para1 = "Working on this project has been a very informative experience into the process of agile development. Overall, the project ran smoothly other than a setback at the beginning due to overestimation of the time we had available to us. This gave us a valuable insight into agile development and the sprint process. In future, we would look to fully recreate the site with our own creative twist while making the site mobile friendly and sticking to client requirements. "
$para3 = "We would also look at fleshing out the accessibility features with more options to create as comfortable of a user experience as we can. Further, we would develop a fully secure back-end system including a database to store data, such as user preferences and favourites, anonymously.  We would also look at adding a chatbot and possibly a forum to the site to allow interaction between users. "

$tr20.Text = $para1
[void]$tr20.InsertAfter([char]13)
[void]$tr20.InsertAfter([char]13)
[void]$tr20.InsertAfter($para3)

# The shape auto-fits to its text (<a:spAutoFit/>); pin the height to the
# value PowerPoint computed for the new three-paragraph body.
$shape20.Height = 585.1673279346456

# ---------------------------------------------------------------
# Shape 5 ("TextBox 23") - "We worked well as a group..." paragraph
# ---------------------------------------------------------------
$shape23 = $s.Shapes.Item("TextBox 23")
$tr23 = $shape23.TextFrame.TextRange

$target = "We worked well as a group by effectively delegating tasks to group members and making sure every group member played their part. Group members from the computing science background worked on the front-end prototype and coding, while the cyber security researched security "
$target += "guidelines"
$target += " that our site would have to follow and made sure that they were implemented into the code. Everyone in the group checked on each other" + [char]8217 + "s progress to make sure that we progressed at a speed that would allow us to have a finished final product.  "

$startIdx = 241
$len23 = $target.Length
$sub23 = $tr23.Characters($startIdx, $len23)
$sub23.Text = $target

# This shape's on-slide size does not change in the target edit, so restore
# the pre-edit (auto-fit) height after the text rewrite (value taken from
# the shape's own original <a:ext cy="6994353"/> converted to points).
$shape23.Height = 550.7364807929134
